$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "'2024-03-04"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "Erbseneintopf i Baguette a,a-1 Wiener Würstchen 1,2,16"
$ws.Cells.Item(2, 3).Value = "Erbseneintopf i Baguette a,a-1"

$ws.Cells.Item(3, 1).Value = "'2024-03-05"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "Currywurst  Currywurstsoße 1  Mayonaise-Kartoffelsalat c,k,12"
$ws.Cells.Item(3, 3).Value = "Tofuwurst (Curry Brat) a,c,a-1 Currywurstsoße 1  Mayonaise-Kartoffelsalat c,k,12"

$ws.Cells.Item(4, 1).Value = "'2024-03-06"
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(4, 2).Value = "Frikadelle a,c,a-1 Kartoffeln  Champignon-Gemüseragout g,p"
$ws.Cells.Item(4, 3).Value = "Champignonpfanne g,p Semmelknödel a,c,a-1"

$ws.Cells.Item(5, 1).Value = "'2024-03-07"
$ws.Cells.Item(5, 1).Style = "Normal"
$ws.Cells.Item(5, 2).Value = "Hühnerfrikassee g,p Reis"
$ws.Cells.Item(5, 3).Value = "Bohnen-Gemüseeintopf to- matisiert  Baguette a,a-1"

$ws.Cells.Item(6, 1).Value = "'2024-03-08"
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(6, 2).Value = "Fischfilet d Kräutersoße g,p Kartoffeln  Blattspinat"
$ws.Cells.Item(6, 3).Value = "Vegetarische Maultaschen a,c,g,i,p,a-1 Kräutersoße g,p Blattspinat"

$ws.Cells.Item(7, 1).Value = "'2024-03-25"
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).Value = "Chili con Carne (Rind) 1  Baguette a,a-1"
$ws.Cells.Item(7, 3).Value = "Chili sin Carne (Vegetarisch) a,f,a-1,a-5,1  Baguette a,a-1"

$ws.Cells.Item(8, 1).Value = "'2024-03-26"
$ws.Cells.Item(8, 1).Style = "Normal"
$ws.Cells.Item(8, 2).Value = "Fleischkäse 1,2,16  Süßer Senf k Röstkartoffeln 2  Bayrisch Kraut"
$ws.Cells.Item(8, 3).Value = "Gemüsefrikadelle a,b,c,e,f,g,h,i,k,l,m,p,a-1 Röstkartoffeln 2  Bayrisch Kraut"

$ws.Cells.Item(9, 1).Value = "'2024-03-27"
$ws.Cells.Item(9, 1).Style = "Normal"
$ws.Cells.Item(9, 2).Value = "Putenkeule geschmort a,c,f,g,i,k,p,2  Champignon-Gemüsesoße g,p Kräuterspätzle a,c,a-1"
$ws.Cells.Item(9, 3).Value = "Kräuterspätzle a,c,a-1 Champignon-Gemüseragout g,p"

$ws.Cells.Item(10, 1).Value = "'2024-03-28"
$ws.Cells.Item(10, 1).Style = "Normal"
$ws.Cells.Item(10, 2).Value = "Bauernfrühstück c,g,p,1,2,16  Gewürzgurke 1,4  Salatbeilage c,g,k,p"
$ws.Cells.Item(10, 3).Value = "Bauernfrühstück VEG c,g,p Gewürzgurke 1,4  Salatbeilage c,g,k,p"

$ws.Cells.Item(11, 1).Value = "'2024-03-29"
$ws.Cells.Item(11, 1).Style = "Normal"
$ws.Cells.Item(11, 2).Value = "Fischstäbchen a,d,a-1 Kartoffeln  Rahmspinat g,p"
$ws.Cells.Item(11, 3).Value = "Gemüsestäbchen gebraten a,c,g,i,p,a-1 Kartoffeln  Rahmspinat g,p"

$ws.Cells.Item(12, 1).Value = "'2024-03-11"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = "Karotten-Weißkohleintopf g,p Baguette a,a-1 Kassler Einlage 1,2"
$ws.Cells.Item(12, 3).Value = "Karotten-Weißkohleintopf g,p Baguette a,a-1"

$ws.Cells.Item(13, 1).Value = "'2024-03-12"
$ws.Cells.Item(13, 1).Style = "Normal"
$ws.Cells.Item(13, 2).Value = "Penne a,a-1,2  Arrabiatasoße mit Speck 1,2  Salatbeilage c,g,k,p"
$ws.Cells.Item(13, 3).Value = "Penne a,a-1,2  Arrabiatasoße mit Mozzarel- la g,p,1  Salatbeilage c,g,k,p"

$ws.Cells.Item(14, 1).Value = "'2024-03-13"
$ws.Cells.Item(14, 1).Style = "Normal"
$ws.Cells.Item(14, 2).Value = "Hähnchenstreifen Zürcher Art g,p Spätzle a,c,a-1"
$ws.Cells.Item(14, 3).Value = "Tofugeschnetzeltes Zürcher Art f,g,p Spätzle a,c,a-1"

$ws.Cells.Item(15, 1).Value = "'2024-03-14"
$ws.Cells.Item(15, 1).Style = "Normal"
$ws.Cells.Item(15, 2).Value = "Kohlwurst k,1,2,3,16  Senf k Röstkartoffeln 2  Grünkohl k,1,2,3,16"
$ws.Cells.Item(15, 3).Value = "Spanische Gemüsepfanne  Hirtenkäse (30g) g,p Röstkartoffeln 2  Sour Cream g,p,2"

$ws.Cells.Item(16, 1).Value = "'2024-03-15"
$ws.Cells.Item(16, 1).Style = "Normal"
$ws.Cells.Item(16, 2).Value = "Fischfilet d Finkenwerder Garnitur 1,2  Kartoffeln  Gurkensalat in Dill-Sahne- Joghurt-Dressing g,p"
$ws.Cells.Item(16, 3).Value = "Paprika-Kartoffelcurry g,i,p,1  Koriander-Dip g,p"

$ws.Cells.Item(17, 1).Value = "'2024-03-18"
$ws.Cells.Item(17, 1).Style = "Normal"
$ws.Cells.Item(17, 2).Value = "Tomatisierter Gemüseeintopf  Baguette a,a-1 Chorizo f,g,p,1,2,12,16"
$ws.Cells.Item(17, 3).Value = "Tomatisierter Gemüseeintopf  Baguette a,a-1"

$ws.Cells.Item(18, 1).Value = "'2024-03-19"
$ws.Cells.Item(18, 1).Style = "Normal"
$ws.Cells.Item(18, 2).Value = "Kohlroulade g,k,p Kümmel-Specksoße 1,2  Kartoffeln"
$ws.Cells.Item(18, 3).Value = "Kohlroulade VEG a,c,f,a-1,a-2 Kümmel-Kräutersoße  Kartoffeln"

$ws.Cells.Item(19, 1).Value = "'2024-03-20"
$ws.Cells.Item(19, 1).Style = "Normal"
$ws.Cells.Item(19, 2).Value = "Hähnchenschnitte gefüllt a,g,i,p,a-1 Kartoffelpüree g,m,p,2  gestovter Rosenkohl g,p"
$ws.Cells.Item(19, 3).Value = "Brokkoli-Kartoffelauflauf g,p"

$ws.Cells.Item(20, 1).Value = "'2024-03-21"
$ws.Cells.Item(20, 1).Style = "Normal"
$ws.Cells.Item(20, 2).Value = "Schweinebraten 2  Bratensoße  Röstkartoffeln 2  Kohlrabi g,p"
$ws.Cells.Item(20, 3).Value = "Roter Linsenbratling a,b,c,e,f,g,h,i,k,l,m,p Bratensoße  Röstkartoffeln 2  Kohlrabi g,p"

$ws.Cells.Item(21, 1).Value = "'2024-03-22"
$ws.Cells.Item(21, 1).Style = "Normal"
$ws.Cells.Item(21, 2).Value = "Überbackenes Schlemmerfi- let a,c,d,g,p,a-1 Kartoffeln  Wurzelgemüse gestovt g,i,p"
$ws.Cells.Item(21, 3).Value = "Eieromelette c,g,p Kartoffeln  Wurzelgemüse gestovt g,i,p"
